$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2704
$ws1.Range("F3").Value = 587
$ws1.Range("F4").Value = 477
$ws1.Range("F6").Value = 211
$ws1.Range("F7").Value = 508
$ws1.Range("F8").Value = 1262
$ws1.Range("F9").Value = 592
$ws1.Range("F13").Value = 381
$ws1.Range("F14").Value = 5895
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 1830
$ws1.Range("F17").Value = 4337
$ws1.Range("F21").Value = 5055
$ws1.Range("F22").Value = 6520
$ws1.Range("F25").Value = 713
$ws1.Range("F26").Value = 3856
$ws1.Range("F27").Value = 516
$ws1.Range("F29").Value = 207
$ws1.Range("F30").Value = 135
$ws1.Range("F31").Value = 1009
$ws1.Range("F32").Value = 1440
$ws1.Range("F34").Value = 611
$ws1.Range("F35").Value = 1632
$ws1.Range("F37").Value = 1778
$ws1.Range("F38").Value = 213
$ws1.Range("F39").Value = 1172
$ws1.Range("F41").Value = 646
$ws1.Range("F43").Value = 3518
$ws1.Range("F45").Value = 312
$ws1.Range("F46").Value = 422
$ws1.Range("F47").Value = 12
$ws1.Range("F48").Value = 56
$ws1.Range("F49").Value = 3909

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7
$ws2.Range("F4").Value = 1221
$ws2.Range("F10").Value = 18
$ws2.Range("F11").Value = 18
$ws2.Range("F28").Value = 49

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 4109

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4109
$ws4.Range("F3").Value = 2704
$ws4.Range("F4").Value = 587
$ws4.Range("F5").Value = 477
$ws4.Range("F7").Value = 1221
$ws4.Range("F9").Value = 211
$ws4.Range("F10").Value = 508
$ws4.Range("F12").Value = 1262
$ws4.Range("F13").Value = 18
$ws4.Range("F14").Value = 592
$ws4.Range("F17").Value = 381
$ws4.Range("F18").Value = 104
$ws4.Range("F19").Value = 1830
$ws4.Range("F20").Value = 4338
$ws4.Range("F21").Value = 5055
$ws4.Range("F22").Value = 5055
$ws4.Range("F25").Value = 713
$ws4.Range("F26").Value = 3856
$ws4.Range("F27").Value = 516
$ws4.Range("F28").Value = 207
$ws4.Range("F29").Value = 135
$ws4.Range("F30").Value = 1009
$ws4.Range("F31").Value = 1440
$ws4.Range("F33").Value = 611
$ws4.Range("F34").Value = 1632
$ws4.Range("F36").Value = 1778
$ws4.Range("F40").Value = 646
$ws4.Range("F44").Value = 3518
$ws4.Range("F47").Value = 312
$ws4.Range("F48").Value = 56
$ws4.Range("F50").Value = 3909
